# The "loop" (e.g. a gradient-descent / matrix-factorization fitting
# routine) that produces column C ("predicted factors") was re-run, and
# as the commit message notes, running it again yields slightly
# different numbers than before. Update column C (rows 2-8) with the
# newly computed predicted values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = -4.969492719029112
    3 = -1.160214231626897
    4 = -0.05269727303333466
    5 = -0.3898423783531488
    6 = 0.01194139783248615
    7 = 0.1057079722102457
    8 = 0.1289825989495366
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}
